$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Change Runmode column (D) from Y to N for existing test rows 2-45 ---
$ws.Range("D2:D45").Value = "N"

# --- 2. Row 45 Results (E) changes from PASS to SKIP ---
$ws.Range("E45").Value = "SKIP"

# --- 3. Build new rows 46 and 47 ---
# Use row 27 as a formatting template: it already carries the border/fill/wrap
# styles (A=s7, B=s8, C=s4, D=s7, E=s3) and is a wrapped, multi-line row.
$ws.Range("A27:E27").Copy($ws.Range("A46:E46"))
$ws.Range("A27:E27").Copy($ws.Range("A47:E47"))

# Row 46: TestCase_B45
$ws.Range("A46").Value = "TestCase_B45"
$ws.Range("B46").Value = "OPQA-270"
$ws.Range("C46").Value = "Verify that following  content type options are present in the search drop down:`na)All`nb)Articles`nc)Patents`nd)People`ne)Posts"
$ws.Range("D46").Value = "N"
$ws.Range("E46").Value = "SKIP"
$ws.Rows.Item(46).RowHeight = 90

# Row 47: TestCase_B46
$ws.Range("A47").Value = "TestCase_B46"
$ws.Range("B47").Value = "OPQA-274"
$ws.Range("C47").Value = "Verify that ALL content type is selected in the search drop down by default"
$ws.Range("D47").Value = "Y"
$ws.Range("E47").Value = "PASS"

# D47 uses a plain bordered style (s=3) rather than the filled one (s=7) used
# elsewhere in the D column - align it with the style already used on column E.
$ws.Range("E2").Copy($ws.Range("D47"))
$ws.Range("D47").Value = "Y"

# --- 4. Update sheet view / selection to match the edited range ---
$ws.Range("D48").Select() | Out-Null
